# Insert a new row at row 200 (shifts existing rows 200-290 down to 201-291)
# and populate it with a new data record, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(200).Insert()

$ws.Cells.Item(200, 1).Value2 = 8
$ws.Cells.Item(200, 2).Value2 = "Terminal La Palmera de La Serena"
$ws.Cells.Item(200, 3).Value2 = "Coquimbo"
$ws.Cells.Item(200, 4).Value2 = 44875
$ws.Cells.Item(200, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(200, 5).Value2 = 4
$ws.Cells.Item(200, 6).Value2 = 100112031
$ws.Cells.Item(200, 7).Value2 = "Poroto verde"
$ws.Cells.Item(200, 8).Value2 = "Magnum"
$ws.Cells.Item(200, 9).Value2 = "Primera"
$ws.Cells.Item(200, 10).Value2 = 500
$ws.Cells.Item(200, 11).Value2 = 41000
$ws.Cells.Item(200, 12).Value2 = 42000
$ws.Cells.Item(200, 13).Value2 = 41500
$ws.Cells.Item(200, 14).Value2 = "$/malla 25 kilos"
$ws.Cells.Item(200, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(200, 16).Value2 = 1660
$ws.Cells.Item(200, 17).Value2 = 25
$ws.Cells.Item(200, 18).Value2 = "Hortaliza"
